$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("存款")

# --- Header row (row 1) ---
$ws.Cells.Item(1,2).Value = "bank"
$ws.Cells.Item(1,3).Value = "deposit_type"
$ws.Cells.Item(1,4).Value = "currency"
$ws.Cells.Item(1,5).Value = "owner"
$ws.Cells.Item(1,6).Value = "total"
$ws.Cells.Item(1,7).Value = "property_category"
$ws.Cells.Item(1,8).Value = "category"
$ws.Cells.Item(1,9).Value = "date"
$ws.Cells.Item(1,10).Value = "legislator_name"
$ws.Cells.Item(1,11).Value = "legislator_id"
$ws.Cells.Item(1,12).Value = "source_file"
$ws.Cells.Item(1,13).Value = "index"

# --- Data rows (rows 2-9): bank / deposit_type / currency / owner ---
$banks = @("玉山商業銀行北新分行","玉山商業銀行北新分行","中國信託商業銀行板橋分行","華南商業銀行華江分行","華南商業銀行文山分行","華南商業銀行文山分行","玉山商業銀行北新分行","中國信託商業銀行板橋分行")
$owners = @("羅明才","黃靜秋","黃靜秋","黃靜秋","羅〇偉","羅〇立","黃靜秋","黃靜秋")

for ($i = 0; $i -lt 8; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r,2).Value = $banks[$i]
    $ws.Cells.Item($r,3).Value = "活期儲蓄存款"
    $ws.Cells.Item($r,4).Value = "新臺幣"
    $ws.Cells.Item($r,5).Value = $owners[$i]

    $ws.Cells.Item($r,7).Value = "deposit"
    $ws.Cells.Item($r,8).Value = "normal"
    $ws.Cells.Item($r,9).Value = "2012-04-30"
    $ws.Cells.Item($r,10).Value = "羅明才"
    $ws.Cells.Item($r,11).Value = 879
    $ws.Cells.Item($r,12).Value = "tmpa5201"
    $ws.Cells.Item($r,13).Value = $ws.Cells.Item($r,1).Value
}
